# Apply numeric updates to the Leve profit-tracking tables (H:N columns)
# across all 8 sheets, per the scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 454.8
$ws.Range("I2").Value = 481
$ws.Range("K2").Value = 481
$ws.Range("M2").Value = -368
$ws.Range("H4").Value = 2588.818
$ws.Range("I4").Value = 2036.3334
$ws.Range("K4").Value = 2036.3334
$ws.Range("M4").Value = -1922.3334
$ws.Range("H32").Value = 6332.1665
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 8498.25
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 8498.25
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -9150.25
$ws.Range("H40").Value = 3089.72
$ws.Range("I40").Value = 2499.889
$ws.Range("J40").Value = 3421.5
$ws.Range("K40").Value = 2499.889
$ws.Range("L40").Value = 3421.5
$ws.Range("M40").Value = -2324.889
$ws.Range("N40").Value = -3771.5
$ws.Range("H106").Value = 42250.2
$ws.Range("I106").Value = 2437.6191
$ws.Range("K106").Value = 2437.6191
$ws.Range("M106").Value = -1806.6191
$ws.Range("H107").Value = 319.3
$ws.Range("I107").Value = 319
$ws.Range("K107").Value = 319
$ws.Range("M107").Value = 1601
$ws.Range("H127").Value = 1043.0769
$ws.Range("I127").Value = 395.55554
$ws.Range("K127").Value = 1186.66662
$ws.Range("M127").Value = 3773.33338

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1390240.8
$ws.Range("I45").Value = 2778556.5
$ws.Range("J45").Value = 1925
$ws.Range("K45").Value = 2778556.5
$ws.Range("L45").Value = 1925
$ws.Range("M45").Value = -2778179.5
$ws.Range("N45").Value = -2679
$ws.Range("H46").Value = 6910.25
$ws.Range("J46").Value = 6910.25
$ws.Range("L46").Value = 6910.25
$ws.Range("N46").Value = -7548.25
$ws.Range("H61").Value = 3544.0667
$ws.Range("I61").Value = 2563.5
$ws.Range("K61").Value = 2563.5
$ws.Range("M61").Value = -2351.5
$ws.Range("H74").Value = 37040190
$ws.Range("I74").Value = 58825990
$ws.Range("K74").Value = 58825990
$ws.Range("M74").Value = -58825116
$ws.Range("H77").Value = 37040190
$ws.Range("I77").Value = 58825990
$ws.Range("K77").Value = 294129950
$ws.Range("M77").Value = -294125582
$ws.Range("H102").Value = 6915857.5
$ws.Range("I102").Value = 7409385
$ws.Range("K102").Value = 7409385
$ws.Range("M102").Value = -7407763
$ws.Range("H136").Value = 3544.0667
$ws.Range("I136").Value = 2563.5
$ws.Range("K136").Value = 7690.5
$ws.Range("M136").Value = -5140.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8643.9375
$ws.Range("I20").Value = 2765.3635
$ws.Range("J20").Value = 21576.8
$ws.Range("K20").Value = 2765.3635
$ws.Range("L20").Value = 21576.8
$ws.Range("M20").Value = -2518.3635
$ws.Range("N20").Value = -22070.8
$ws.Range("H105").Value = 2150.718
$ws.Range("I105").Value = 1996.3939
$ws.Range("K105").Value = 1996.3939
$ws.Range("M105").Value = -249.3939
$ws.Range("H107").Value = 55558668
$ws.Range("I107").Value = 2799.8
$ws.Range("J107").Value = 125003500
$ws.Range("K107").Value = 2799.8
$ws.Range("L107").Value = 125003500
$ws.Range("M107").Value = -879.8000000000002
$ws.Range("N107").Value = -125007340

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3318.2
$ws.Range("I31").Value = 2390.3
$ws.Range("J31").Value = 5174
$ws.Range("K31").Value = 2390.3
$ws.Range("L31").Value = 5174
$ws.Range("M31").Value = -2095.3
$ws.Range("N31").Value = -5764
$ws.Range("H34").Value = 3318.2
$ws.Range("I34").Value = 2390.3
$ws.Range("J34").Value = 5174
$ws.Range("K34").Value = 2390.3
$ws.Range("L34").Value = 5174
$ws.Range("M34").Value = -2188.3
$ws.Range("N34").Value = -5578
$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630
$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184
$ws.Range("H105").Value = 902
$ws.Range("I105").Value = 903.1818
$ws.Range("J105").Value = 897.6667
$ws.Range("K105").Value = 903.1818
$ws.Range("L105").Value = 897.6667
$ws.Range("M105").Value = 843.8182
$ws.Range("N105").Value = -4391.6667
$ws.Range("H107").Value = 696.7
$ws.Range("I107").Value = 791.45
$ws.Range("J107").Value = 507.2
$ws.Range("K107").Value = 791.45
$ws.Range("L107").Value = 507.2
$ws.Range("M107").Value = 1128.55
$ws.Range("N107").Value = -4347.2
$ws.Range("H134").Value = 4522.973
$ws.Range("I134").Value = 3985.5667
$ws.Range("K134").Value = 11956.7001
$ws.Range("M134").Value = -9421.7001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 3166.6667
$ws.Range("I115").Value = 3000
$ws.Range("J115").Value = 3250
$ws.Range("K115").Value = 9000
$ws.Range("L115").Value = 9750
$ws.Range("M115").Value = -7825
$ws.Range("N115").Value = -12100
$ws.Range("H131").Value = 9368.634
$ws.Range("J131").Value = 13671.842
$ws.Range("L131").Value = 41015.526
$ws.Range("N131").Value = -51095.526

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 604.6818
$ws.Range("I97").Value = 620.6842
$ws.Range("K97").Value = 620.6842
$ws.Range("M97").Value = -124.6842
$ws.Range("H122").Value = 558445.4
$ws.Range("I122").Value = 697056.75
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 2091170.25
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2088720.25
$ws.Range("N122").Value = -16900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 176667
$ws.Range("I4").Value = 255001.5
$ws.Range("J4").Value = 19998
$ws.Range("K4").Value = 255001.5
$ws.Range("L4").Value = 19998
$ws.Range("M4").Value = -254888.5
$ws.Range("N4").Value = -20224
$ws.Range("H28").Value = 176667
$ws.Range("I28").Value = 255001.5
$ws.Range("J28").Value = 19998
$ws.Range("K28").Value = 255001.5
$ws.Range("L28").Value = 19998
$ws.Range("M28").Value = -254769.5
$ws.Range("N28").Value = -20462
$ws.Range("H32").Value = 5548
$ws.Range("I32").Value = 4659.8
$ws.Range("J32").Value = 9989
$ws.Range("K32").Value = 4659.8
$ws.Range("L32").Value = 9989
$ws.Range("M32").Value = -4342.8
$ws.Range("N32").Value = -10623
$ws.Range("H34").Value = 16931.666
$ws.Range("I34").Value = 16931.666
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 16931.666
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -16759.666
$ws.Range("N34").ClearContents()
$ws.Range("H37").Value = 176667
$ws.Range("I37").Value = 255001.5
$ws.Range("J37").Value = 19998
$ws.Range("K37").Value = 255001.5
$ws.Range("L37").Value = 19998
$ws.Range("M37").Value = -254894.5
$ws.Range("N37").Value = -20212
$ws.Range("H61").Value = 1540.2
$ws.Range("I61").Value = 1425.25
$ws.Range("K61").Value = 1425.25
$ws.Range("M61").Value = -1223.25
$ws.Range("H113").Value = 1540.2
$ws.Range("I113").Value = 1425.25
$ws.Range("K113").Value = 1425.25
$ws.Range("M113").Value = 744.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 7741.222
$ws.Range("J45").Value = 8029.6665
$ws.Range("L45").Value = 8029.6665
$ws.Range("N45").Value = -9011.666499999999
$ws.Range("H62").Value = 2982815.5
$ws.Range("I62").Value = 11907262
$ws.Range("K62").Value = 11907262
$ws.Range("M62").Value = -11906638
$ws.Range("H65").Value = 2982815.5
$ws.Range("I65").Value = 11907262
$ws.Range("K65").Value = 59536310
$ws.Range("M65").Value = -59533190
$ws.Range("H133").Value = 29999.092
$ws.Range("J133").Value = 28999
$ws.Range("L133").Value = 28999
$ws.Range("N133").Value = -39119
